$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = "project coordinator"
$ws.Range("B6").Value = "Sam"
